$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (ListObject) backing the data grid - grow it by one row so the
# table ref / autoFilter ref extend from A1:E25 to A1:E26.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Fill in the new row's data (row 26)
$ws.Range("A26").Value = "1731. The Number of Employees Which Report to Each Employee"
$ws.Range("B26").Value = "Easy"
$ws.Range("C26").Value = "Advanced Select and Joins"
$ws.Range("D26").Value = "Use inner join of Employees e1 and Employees e2 on e1.employee_id = e2.reports_to, group by e1.employee_id, e1.name. For the reports count, select count(e2.employee_id)"
$ws.Range("E26").Value = "https://leetcode.com/problems/the-number-of-employees-which-report-to-each-employee/solutions/3865526/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 "

# Match the "Easy" fill/format used on the other Easy rows (e.g. B2) by
# copying that cell's formatting onto the new Difficulty cell.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null

# Turn the new Link cell into a real hyperlink, styled like the others.
# (The displayed text keeps its trailing space like the other rows, but the
# underlying hyperlink address/target itself is the trimmed URL.)
$ws.Hyperlinks.Add($ws.Range("E26"), "https://leetcode.com/problems/the-number-of-employees-which-report-to-each-employee/solutions/3865526/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50") | Out-Null
$ws.Range("E26").Style = "Hyperlink"

$excel.CutCopyMode = $false

# Match the saved cursor/selection position recorded in the workbook.
$ws.Range("E29:E30").Select() | Out-Null
